$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.166.09"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.829.23"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.16"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9993"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2920"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.13"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07667"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "1.833.13"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.956"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6665"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.59"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008978"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "29.126.96"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "2.075.11"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.84"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.367"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.28"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1416"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.519"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.487"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05887"
$ws.Range("E30").Value = "  +5.30%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.100"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.079"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.207"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7318"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.838"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "1.226.20"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.297"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9205"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.10"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "1.978.47"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.08"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5046"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.130"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4025"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("E51").Value = "  +1.68%  "
